$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.01868167416557
$ws.Range("C2").Value = 10.48541145433043
$ws.Range("E2").Value = 23.09815684862075
$ws.Range("F2").Value = 40.66186628316461
$ws.Range("G2").Value = 27.93917481772626
$ws.Range("H2").Value = 13.83206282455533
$ws.Range("J2").Value = 7.843520715174872
$ws.Range("B3").Value = 15.28907220859163
$ws.Range("C3").Value = 9.820664496775212
$ws.Range("E3").Value = 22.82268482753717
$ws.Range("F3").Value = 40.45422839896959
$ws.Range("G3").Value = 27.99899962343071
$ws.Range("H3").Value = 13.92241251542895
$ws.Range("J3").Value = 7.882584870993004
$ws.Range("B4").Value = 14.82472780281939
$ws.Range("C4").Value = 9.387368127217611
$ws.Range("E4").Value = 22.65534676695471
$ws.Range("F4").Value = 40.34092557981001
$ws.Range("G4").Value = 28.05719956661402
$ws.Range("H4").Value = 13.98267214225235
$ws.Range("J4").Value = 7.908082756673471
$ws.Range("B5").Value = 14.63164381326675
$ws.Range("C5").Value = 9.204477715763066
$ws.Range("E5").Value = 22.58767544609957
$ws.Range("F5").Value = 40.29835129393474
$ws.Range("G5").Value = 28.08624258526624
$ws.Range("H5").Value = 14.00842393783528
$ws.Range("J5").Value = 7.918853559549685
$ws.Range("B6").Value = 14.59935739221329
$ws.Range("C6").Value = 9.173727457678595
$ws.Range("E6").Value = 22.57647203865888
$ws.Range("F6").Value = 40.29149989413391
$ws.Range("G6").Value = 28.09138458029287
$ws.Range("H6").Value = 14.01277198887845
$ws.Range("J6").Value = 7.920665004822725
$ws.Range("B7").Value = 14.82213907014601
$ws.Range("C7").Value = 9.384927179631058
$ws.Range("E7").Value = 22.65443193339106
$ws.Range("F7").Value = 40.34033680736561
$ws.Range("G7").Value = 28.05756978710582
$ws.Range("H7").Value = 13.98301460949495
$ws.Range("J7").Value = 7.908226476148461
$ws.Range("B8").Value = 15.77064958674835
$ws.Range("C8").Value = 10.26142441862444
$ws.Range("E8").Value = 23.00284030566868
$ws.Range("F8").Value = 40.58734802135282
$ws.Range("G8").Value = 27.95530603271231
$ws.Range("H8").Value = 13.8622177678417
$ws.Range("J8").Value = 7.856676097636771
$ws.Range("B9").Value = 17.49143836745913
$ws.Range("C9").Value = 11.78089240717212
$ws.Range("E9").Value = 23.69740492248516
$ws.Range("F9").Value = 41.18268608717467
$ws.Range("G9").Value = 27.92788699940058
$ws.Range("H9").Value = 13.66364235055386
$ws.Range("J9").Value = 7.767589564584886
$ws.Range("B10").Value = 18.66037325616834
$ws.Range("C10").Value = 12.77616991548096
$ws.Range("E10").Value = 24.21063094995593
$ws.Range("F10").Value = 41.68512154896077
$ws.Range("G10").Value = 28.01640872249636
$ws.Range("H10").Value = 13.54157731261365
$ws.Range("J10").Value = 7.709461206090215
$ws.Range("B11").Value = 19.16978066586198
$ws.Range("C11").Value = 13.20280595784114
$ws.Range("E11").Value = 24.44390951703856
$ws.Range("F11").Value = 41.92714197303039
$ws.Range("G11").Value = 28.08073035905919
$ws.Range("H11").Value = 13.4913247393187
$ws.Range("J11").Value = 7.684609181701783
$ws.Range("B12").Value = 19.35935651580592
$ws.Range("C12").Value = 13.3606193063929
$ws.Range("E12").Value = 24.53214980827171
$ws.Range("F12").Value = 42.02065995319459
$ws.Range("G12").Value = 28.10856738056047
$ws.Range("H12").Value = 13.47306274700474
$ws.Range("J12").Value = 7.675427391835949
$ws.Range("B13").Value = 19.3186774631914
$ws.Range("C13").Value = 13.32679771421153
$ws.Range("E13").Value = 24.51315110318057
$ws.Range("F13").Value = 42.00043709920149
$ws.Range("G13").Value = 28.10241719788282
$ws.Range("H13").Value = 13.47696151247498
$ws.Range("J13").Value = 7.677394657606898
$ws.Range("B14").Value = 19.18544441275771
$ws.Range("C14").Value = 13.21586433003767
$ws.Range("E14").Value = 24.45117144233807
$ws.Range("F14").Value = 41.93479866475899
$ws.Range("G14").Value = 28.08295070603937
$ws.Range("H14").Value = 13.48980687942273
$ws.Range("J14").Value = 7.68384919593728
$ws.Range("B15").Value = 19.10339921228105
$ws.Range("C15").Value = 13.14742715883109
$ws.Range("E15").Value = 24.4131924000553
$ws.Range("F15").Value = 41.89483475129563
$ws.Range("G15").Value = 28.07148052648617
$ws.Range("H15").Value = 13.49777526479287
$ws.Range("J15").Value = 7.687832638357466
$ws.Range("B16").Value = 18.62662257534182
$ws.Range("C16").Value = 12.7477635968653
$ws.Range("E16").Value = 24.19537580930289
$ws.Range("F16").Value = 41.66957030610286
$ws.Range("G16").Value = 28.0126912094733
$ws.Range("H16").Value = 13.54496835626648
$ws.Range("J16").Value = 7.711117387240453
$ws.Range("B17").Value = 18.32832593770138
$ws.Range("C17").Value = 12.49590065756773
$ws.Range("E17").Value = 24.06165344117542
$ws.Range("F17").Value = 41.53478021629261
$ws.Range("G17").Value = 27.98280486117816
$ws.Range("H17").Value = 13.57527730750611
$ws.Range("J17").Value = 7.725809535328473
$ws.Range("B18").Value = 18.1546555921016
$ws.Range("C18").Value = 12.3485763549886
$ws.Range("E18").Value = 23.98472603520497
$ws.Range("F18").Value = 41.45852366850796
$ws.Range("G18").Value = 27.96787841741376
$ws.Range("H18").Value = 13.59320622137358
$ws.Range("J18").Value = 7.734409785278364
$ws.Range("B19").Value = 18.09549719809935
$ws.Range("C19").Value = 12.29827214083556
$ws.Range("E19").Value = 23.95867951265974
$ws.Range("F19").Value = 41.43292475027921
$ws.Range("G19").Value = 27.96321249492444
$ws.Range("H19").Value = 13.5993615611435
$ws.Range("J19").Value = 7.737347385103273
$ws.Range("B20").Value = 18.3602981333722
$ws.Range("C20").Value = 12.52296635574338
$ws.Range("E20").Value = 24.07589035819327
$ws.Range("F20").Value = 41.54899775246888
$ws.Range("G20").Value = 27.98575187374879
$ws.Range("H20").Value = 13.57199946859161
$ws.Range("J20").Value = 7.724230033630294
$ws.Range("B21").Value = 19.22466924852744
$ws.Range("C21").Value = 13.24854967805869
$ws.Range("E21").Value = 24.46937955921433
$ws.Range("F21").Value = 41.95402805386085
$ws.Range("G21").Value = 28.08857393299887
$ws.Range("H21").Value = 13.4860129804016
$ws.Range("J21").Value = 7.681947119874483
$ws.Range("B22").Value = 19.77016181164334
$ws.Range("C22").Value = 13.7009410335817
$ws.Range("E22").Value = 24.72594945989663
$ws.Range("F22").Value = 42.22960572760346
$ws.Range("G22").Value = 28.17605895524588
$ws.Range("H22").Value = 13.43429401595087
$ws.Range("J22").Value = 7.655648545544993
$ws.Range("B23").Value = 19.48083178748002
$ws.Range("C23").Value = 13.4614835336089
$ws.Range("E23").Value = 24.58909066013836
$ws.Range("F23").Value = 42.08155316070699
$ws.Range("G23").Value = 28.12750642274441
$ws.Range("H23").Value = 13.46148466578133
$ws.Range("J23").Value = 7.669562237379557
$ws.Range("B24").Value = 18.34585027162527
$ws.Range("C24").Value = 12.51073781920966
$ws.Range("E24").Value = 24.06945399419484
$ws.Range("F24").Value = 41.54256614917811
$ws.Range("G24").Value = 27.98441250509249
$ws.Range("H24").Value = 13.57347981061763
$ws.Range("J24").Value = 7.724943648116082
$ws.Range("B25").Value = 17.04201332252383
$ws.Range("C25").Value = 11.391141466106
$ws.Range("E25").Value = 23.5087051375807
$ws.Range("F25").Value = 41.00999815102056
$ws.Range("G25").Value = 27.91640966336871
$ws.Range("H25").Value = 13.71321251958932
$ws.Range("J25").Value = 7.790404278016623
